# planilha_gincana_solidaria.xlsx — "Add files via upload"
#
# The sheet "doacoes_registros" (3rd tab, the active one) had the
# donation-entry rows for the 1st sprint (rows 3 and 4, columns C:H —
# Nome/Grupo related lookups aside: Categoria, Tipo_Item, Quantidade,
# Pontos_Unit, plus the dependent "Novo"/"Brinquedos" lookups) wiped back
# to blank templates, ready for new data entry. The dependent formula
# cells (Pontos_Total = Quantidade*Pontos_Unit, Total_Geral) recalculate
# to 0 once their inputs are gone. Cell formatting/styles are left intact
# — only the values are cleared (Excel's "Clear Contents" / Delete key).

$wb = $excel.ActiveWorkbook

$wsParticipantes = $wb.Worksheets.Item("participantes")
$wsCategorias    = $wb.Worksheets.Item("categorias")
$wsDoacoes       = $wb.Worksheets.Item("doacoes_registros")

# Leave the other two sheets' selection where it was.
$wsParticipantes.Activate()
$wsParticipantes.Range("B10").Select()

$wsCategorias.Activate()
$wsCategorias.Range("A3").Select()

# Do the actual edit on the active sheet: clear C3:H4 (values only,
# formatting/styles untouched) and leave that range selected.
$wsDoacoes.Activate()
$wsDoacoes.Range("C3:H4").ClearContents()
$wsDoacoes.Range("C3:H4").Select()
